# Added jasmine and oak to the voedselbos (food forest) species list.
#
# Sheet1 columns: A=abbr, B=species, C=name_nl, D=name_en, E=height, F=width, G=races
# Two new rows are appended right after the existing data (which runs to row 127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 128: Jasmine.
# (Column D is written before C so that new shared-string entries are
# registered in the same order as the source workbook: abbr, species,
# name_en, name_nl.)
$ws.Cells.Item(128, 1).Value = "J. officinale"
$ws.Cells.Item(128, 2).Value = "Jasminum officinale"
$ws.Cells.Item(128, 4).Value = "Common jasmine"
$ws.Cells.Item(128, 3).Value = "Zomerjasmijn"
$ws.Cells.Item(128, 5).Value = 2
$ws.Cells.Item(128, 6).Value = 2

# Row 129: Oak.
$ws.Cells.Item(129, 1).Value = "Q. robur"
$ws.Cells.Item(129, 2).Value = "Quercus robur"
$ws.Cells.Item(129, 3).Value = "Zomereik"
$ws.Cells.Item(129, 4).Value = "Common oak"
$ws.Cells.Item(129, 5).Value = 6
$ws.Cells.Item(129, 6).Value = 4

# Scroll the view down to the newly added rows and select the cell the
# author landed on after typing the new data.
$excel.ActiveWindow.ScrollRow = 118
$ws.Range("D133").Select()
